# Update NATMI LR-pair edge-weight statistics for Il15-Il2rg (rows 2-10)
# following recomputation with the new ligand/receptor-expressing cell counts
# ("Natmi following Dr Hou advice").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.79549166666667
$ws.Range("H2").Value = 38.386475
$ws.Range("I2").Value = 0.5145949251267348
$ws.Range("J2").Value = 0.5145949251267348
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 32.30682666666667
$ws.Range("N2").Value = 96.92048
$ws.Range("O2").Value = 0.886587237369156
$ws.Range("P2").Value = 0.8865872373691559
$ws.Range("Q2").Value = 413.3817313897777
$ws.Range("R2").Value = 3720.435582508
$ws.Range("S2").Value = 0.4562332930322995
$ws.Range("T2").Value = 0.4562332930322994

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.79549166666667
$ws.Range("H3").Value = 38.386475
$ws.Range("I3").Value = 0.5145949251267348
$ws.Range("J3").Value = 0.5145949251267348
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.7369563333333332
$ws.Range("N3").Value = 2.210869
$ws.Range("O3").Value = 0.02022408719906369
$ws.Range("P3").Value = 0.02022408719906369
$ws.Range("Q3").Value = 9.429718621863886
$ws.Range("R3").Value = 84.86746759677499
$ws.Range("S3").Value = 0.01040721263795873
$ws.Range("T3").Value = 0.01040721263795873

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.79549166666667
$ws.Range("H4").Value = 38.386475
$ws.Range("I4").Value = 0.5145949251267348
$ws.Range("J4").Value = 0.5145949251267348
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.395752
$ws.Range("N4").Value = 10.187256
$ws.Range("O4").Value = 0.09318867543178035
$ws.Range("P4").Value = 0.09318867543178033
$ws.Range("Q4").Value = 43.45031641806666
$ws.Range("R4").Value = 391.0528477626
$ws.Range("S4").Value = 0.0479544194564766
$ws.Range("T4").Value = 0.04795441945647659

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.803896
$ws.Range("H5").Value = 29.411688
$ws.Range("I5").Value = 0.3942822409249843
$ws.Range("J5").Value = 0.3942822409249843
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 32.30682666666667
$ws.Range("N5").Value = 96.92048
$ws.Range("O5").Value = 0.886587237369156
$ws.Range("P5").Value = 0.8865872373691559
$ws.Range("Q5").Value = 316.7327687300266
$ws.Range("R5").Value = 2850.59491857024
$ws.Range("S5").Value = 0.3495656027254018
$ws.Range("T5").Value = 0.3495656027254017

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 9.803896
$ws.Range("H6").Value = 29.411688
$ws.Range("I6").Value = 0.3942822409249843
$ws.Range("J6").Value = 0.3942822409249843
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.7369563333333332
$ws.Range("N6").Value = 2.210869
$ws.Range("O6").Value = 0.02022408719906369
$ws.Range("P6").Value = 0.02022408719906369
$ws.Range("Q6").Value = 7.225043248541332
$ws.Range("R6").Value = 65.025389236872
$ws.Range("S6").Value = 0.007973998421509119
$ws.Range("T6").Value = 0.007973998421509119

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.803896
$ws.Range("H7").Value = 29.411688
$ws.Range("I7").Value = 0.3942822409249843
$ws.Range("J7").Value = 0.3942822409249843
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.395752
$ws.Range("N7").Value = 10.187256
$ws.Range("O7").Value = 0.09318867543178035
$ws.Range("P7").Value = 0.09318867543178033
$ws.Range("Q7").Value = 33.291599449792
$ws.Range("R7").Value = 299.624395048128
$ws.Range("S7").Value = 0.03674263977807338
$ws.Range("T7").Value = 0.03674263977807338

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.265785
$ws.Range("H8").Value = 6.797355
$ws.Range("I8").Value = 0.09112283394828093
$ws.Range("J8").Value = 0.09112283394828093
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 32.30682666666667
$ws.Range("N8").Value = 96.92048
$ws.Range("O8").Value = 0.886587237369156
$ws.Range("P8").Value = 0.8865872373691559
$ws.Range("Q8").Value = 73.20032325893334
$ws.Range("R8").Value = 658.8029093304
$ws.Range("S8").Value = 0.08078834161145473
$ws.Range("T8").Value = 0.08078834161145472

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.265785
$ws.Range("H9").Value = 6.797355
$ws.Range("I9").Value = 0.09112283394828093
$ws.Range("J9").Value = 0.09112283394828093
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.7369563333333332
$ws.Range("N9").Value = 2.210869
$ws.Range("O9").Value = 0.02022408719906369
$ws.Range("P9").Value = 0.02022408719906369
$ws.Range("Q9").Value = 1.669784605721667
$ws.Range("R9").Value = 15.028061451495
$ws.Range("S9").Value = 0.001842876139595834
$ws.Range("T9").Value = 0.001842876139595834

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.265785
$ws.Range("H10").Value = 6.797355
$ws.Range("I10").Value = 0.09112283394828093
$ws.Range("J10").Value = 0.09112283394828093
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.395752
$ws.Range("N10").Value = 10.187256
$ws.Range("O10").Value = 0.09318867543178035
$ws.Range("P10").Value = 0.09318867543178033
$ws.Range("Q10").Value = 7.694043945320001
$ws.Range("R10").Value = 69.24639550788
$ws.Range("S10").Value = 0.008491616197230366
$ws.Range("T10").Value = 0.008491616197230366

